$d = $word.ActiveDocument

# The resume currently starts with a centered paragraph holding the
# candidate's name ("Dheeraj Chand"). A new, centered paragraph with the
# contact information needs to be inserted right after it (and before the
# "PROFESSIONAL SUMMARY" heading).
$nameParagraph = $d.Paragraphs.Item(1)
$nameRange = $nameParagraph.Range

# Create a brand-new paragraph mark immediately after the name paragraph.
$nameRange.InsertParagraphAfter()

# That new paragraph is now the second paragraph in the document. Populate
# it with the contact info as plain (unformatted) text, centered, using an
# OOXML fragment so it gets no inherited character formatting (no bold /
# font-size) from the name run before it.
$contactParagraph = $d.Paragraphs.Item(2)
$contactRange = $contactParagraph.Range

$contactXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:t>202.550.7110 | dheeraj.chand@gmail.com | https://www.dheerajchand.com | https://www.linkedin.com/in/dheerajchand/ | Austin, TX</w:t></w:r></w:p>'
[void]$contactRange.InsertXML($contactXml)
